# Generate Report for Archive
# - Update localization status text from "Ready for handoff" to "In Translation"
#   on all three sheets (Overview uses one column per language; the per-language
#   sheets use a single Status column).
# - Re-fit ("AutoFit"-style) the now-narrower Status columns so they match the
#   shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
